$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "Analysis of Ecological Data: Exploratory and Euclidean Methods in Environmental Sciences (ade4)"
$ws.Range("B66").Value = "CRAN"
$ws.Range("C66").Value = 2023
$ws.Range("D66").Value = "Dray, Dufour, Thioulouse & Siberchicot"
$ws.Range("E66").Value = "Software/Package"
$ws.Range("F66").Value = "NA"
$ws.Range("G66").Value = "ade4 package description"
$ws.Range("H66").Value = "Multivariate statistics"
$ws.Range("I66").Value = "yes"
$ws.Range("J66").Value = "yes"
$ws.Range("K66").Value = "yes"

$ws.Range("A66:K66").Select()
